$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.486.78"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "3.814.93"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'670.73"
$ws.Range("E5").Value = "  +7.22%  "
$ws.Range("D6").Value = "'169.81"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("D7").Value = "3.813.32"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("E11").Value = "  +6.72%  "
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'35.84"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "4.453.10"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "3.812.87"
$ws.Range("E16").Value = "  +4.67%  "
$ws.Range("D17").Value = "70.479.70"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "'11.78"
$ws.Range("E18").Value = "  +23.21%  "
$ws.Range("D19").Value = "'17.65"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'7.19"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "'476.23"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").Value = "'0.713"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").Value = "'83.43"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("E25").Value = "  -4.20%  "
$ws.Range("D26").Value = "'12.23"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").Value = "'10.29"
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "3.965.97"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "'2.86"
$ws.Range("E31").Value = "  +7.51%  "
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").Value = "'7.38"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.771.08"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'9.09"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("D42").Value = "'0.965"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "'2.10"
$ws.Range("E44").Value = "  +10.86%  "
$ws.Range("D46").Value = "'45.80"
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("D47").Value = "'159.03"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("D48").Value = "'48.11"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("D51").Value = "'0.000291"
$ws.Range("E51").Value = "  +5.68%  "
